# Update "want to go" counts (column F) for a handful of rows across
# three worksheets, reflecting refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 1281
$ws1.Range("F16").Value = 17
$ws1.Range("F18").Value = 1592
$ws1.Range("F27").Value = 2741
$ws1.Range("F28").Value = 1516
$ws1.Range("F31").Value = 536
$ws1.Range("F33").Value = 1492
$ws1.Range("F35").Value = 1555
$ws1.Range("F36").Value = 181
$ws1.Range("F38").Value = 814

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 183

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 183
$ws4.Range("F9").Value  = 1281
$ws4.Range("F20").Value = 1592
$ws4.Range("F27").Value = 2741
$ws4.Range("F29").Value = 1516
$ws4.Range("F34").Value = 536
$ws4.Range("F36").Value = 1492
$ws4.Range("F40").Value = 1555
$ws4.Range("F41").Value = 814
